$wb = $excel.ActiveWorkbook

# Add the new "fields" worksheet after the existing "login" sheet
$loginSheet = $wb.Worksheets.Item("login")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$ws.Name = "fields"

# --- Write string-valued cells first, in the exact order they were
# originally authored, so shared-string table indices line up with the
# target workbook (header row, then column D, then column C, then column B).
$ws.Range("B1").Value = "fname"
$ws.Range("C1").Value = "lname"
$ws.Range("D1").Value = "post"

$ws.Range("D5").Value = "dfdfre"
$ws.Range("D6").Value = "e34re2"
$ws.Range("D7").Value = "aaaaaa"

$ws.Range("C2").Value = "fdfsfbf"
$ws.Range("C3").Value = "ererer"
$ws.Range("C6").Value = "rrgds"

$ws.Range("B2").Value = "sdvsbsfb"
$ws.Range("B4").Value = "qddddd"
$ws.Range("B5").Value = "shane"
$ws.Range("B7").Value = "fgvdse332"

# A1 reuses the existing "index" shared string from the login sheet.
$ws.Range("A1").Value = "index"

# --- Numeric cells ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

$ws.Range("B3").Value = 333453
$ws.Range("B6").Value = 313

$ws.Range("C4").Value = 343434
$ws.Range("C7").Value = 33333

$ws.Range("D2").Value = 463723
$ws.Range("D3").Value = 342113
$ws.Range("D4").Value = 403601

# Set view: selection on fields sheet
$ws.Range("M6").Select()

# Make "fields" the active (selected) tab
$ws.Activate()
